$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "90.872.41"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").Value = "3.161.51"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.73"
$ws.Range("E5").Value = "  -0.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "626.43"
$ws.Range("E6").Value = "  +1.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.15"
$ws.Range("E7").Value = "  +26.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.370"
$ws.Range("E8").Value = "  -2.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.07%  "

# Row 10
$ws.Range("D10").Value = "3.159.27"
$ws.Range("E10").Value = "  +1.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.748"
$ws.Range("E11").Value = "  +10.99%  "

# Row 12
$ws.Range("E12").Value = "  +6.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.72"
$ws.Range("E13").Value = "  +6.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000247"
$ws.Range("E14").Value = "  -3.46%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "35.32"
$ws.Range("E15").Value = "  +6.86%  "

# Row 16
$ws.Range("D16").Value = "90.647.30"
$ws.Range("E16").Value = "  -0.40%  "

# Row 17
$ws.Range("D17").Value = "3.739.63"
$ws.Range("E17").Value = "  +2.28%  "

# Row 18
$ws.Range("D18").Value = "3.168.34"
$ws.Range("E18").Value = "  +2.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.75"
$ws.Range("E19").Value = "  +5.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.65"
$ws.Range("E20").Value = "  +6.17%  "

# Row 21
$ws.Range("B21").Value = "PEPE"
$ws.Range("C21").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000215"
$ws.Range("E21").Value = "  -4.06%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.49"
$ws.Range("E22").Value = "  +8.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.22"
$ws.Range("E23").Value = "  +8.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.18"
$ws.Range("E24").Value = "  +1.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.56"
$ws.Range("E25").Value = "  +14.84%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.93"
$ws.Range("E26").Value = "  +6.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.41"
$ws.Range("E27").Value = "  +4.75%  "

# Row 28
$ws.Range("D28").Value = "3.320.35"
$ws.Range("E28").Value = "  +2.20%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.165"
$ws.Range("E30").Value = "  -2.18%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.221"
$ws.Range("E31").Value = "  +55.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.26"
$ws.Range("E32").Value = "  +6.30%  "

# Row 33
$ws.Range("E33").Value = "  -1.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.06"
$ws.Range("E34").Value = "  +17.41%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "522.81"
$ws.Range("E35").Value = "  +0.56%  "

# Row 36
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.95"
$ws.Range("E36").Value = "  +5.87%  "

# Row 37
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.66"
$ws.Range("E37").Value = "  -6.91%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.07"
$ws.Range("E38").Value = "  +0.48%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.146"
$ws.Range("E39").Value = "  +4.36%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.33"
$ws.Range("E40").Value = "  +3.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0912"
$ws.Range("E41").Value = "  +25.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.432"
$ws.Range("E42").Value = "  +16.36%  "

# Row 43
$ws.Range("E43").Value = "  -0.28%  "

# Row 44
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("E45").Value = "  +6.11%  "

# Row 46
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.749"
$ws.Range("E46").Value = "  +22.86%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.76"
$ws.Range("E48").Value = "  +13.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "150.84"
$ws.Range("E49").Value = "  +5.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.38"
$ws.Range("E50").Value = "  +10.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.28"
$ws.Range("E51").Value = "  +3.93%  "
